$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.24
$ws.Range("D3").Value = 0.25
$ws.Range("D4").Value = 0.23
$ws.Range("D5").Value = 0.17
$ws.Range("D6").Value = 0.19
$ws.Range("D7").Value = 0.13
$ws.Range("D8").Value = 0.06
